# Update "想去人数" (interested-count) figures to the values captured
# for the newly generated gh-pages output (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 1056   # 苏州·创世次元动漫游戏嘉年华... : 1052 -> 1056
$wsExpo.Range("F8").Value = 204    # 苏州·漫语堂动漫嘉年华 : 203 -> 204
$wsExpo.Range("F9").Value = 381    # 苏州·第三届华盟国漫次元嘉年华 : 380 -> 381
$wsExpo.Range("F15").Value = 12433 # 苏州·COME IN JOY 动漫品牌国潮文化节 : 12417 -> 12433
$wsExpo.Range("F16").Value = 120   # 苏州·Good jump ACG元旦跨年盛典国潮文化节 : 116 -> 120
$wsExpo.Range("F17").Value = 5491  # 苏州·星部落&青铜树动漫嘉年华 : 5487 -> 5491

# --- Sheet "演出" (performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 119    # 苏州·乐队番同人only live Band Set... : 118 -> 119

# --- Sheet "全部类型" (all types, combined list) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 119     # 苏州·乐队番同人only live Band Set... : 118 -> 119
$wsAll.Range("F7").Value = 1056    # 苏州·创世次元动漫游戏嘉年华... : 1052 -> 1056
$wsAll.Range("F10").Value = 204    # 苏州·漫语堂动漫嘉年华 : 203 -> 204
$wsAll.Range("F11").Value = 381    # 苏州·第三届华盟国漫次元嘉年华 : 380 -> 381
$wsAll.Range("F17").Value = 12433  # 苏州·COME IN JOY 动漫品牌国潮文化节 : 12417 -> 12433
$wsAll.Range("F19").Value = 120    # 苏州·Good jump ACG元旦跨年盛典国潮文化节 : 116 -> 120
$wsAll.Range("F20").Value = 5491   # 苏州·星部落&青铜树动漫嘉年华 : 5487 -> 5491
